$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with data that used to be in a later row (9530 DOI)
$ws.Range("A2").Value = 9530
$ws.Range("B2").Value = "10.5327/1516-3180.142s1.12066"
$ws.Range("C2").Value = "https://www.apm.org.br/wp-content/uploads/SPMJ_v142Suppl1.pdf"

# Update row 3 with data that used to be in row 8 (9701 DOI)
$ws.Range("A3").Value = 9701
$ws.Range("B3").Value = "10.5327/1516-3180.142s1.12677"
$ws.Range("C3").Value = "https://www.apm.org.br/wp-content/uploads/SPMJ_v142Suppl1.pdf"

# Delete rows 4 through 8 (old duplicate rows)
$ws.Range("A4:F8").Delete()
